$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add P1 and Q1, matching style of O1 (copies style + value, then overwrite value)
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# New data values for rows 2-25, columns B..O (14 cols)
$data = @(
  @(3.434417643568509,1.144712818294295,0.3123760442629475,0,1.205266753805645,0.73146128074724,0.0003601650289539737,0.001077091442577949,0.4948503304475622,0.3168370846712563,0.5124944197761181,0,0,0),
  @(2.990033006693295,0.9968553770831079,0.2755548077891774,0,1.111939355810662,0.6746673705668655,[double]"7.625252296250551E-08",0.0004708235328751975,0.4743501209024004,0.3172502880096353,0.4507008927300262,0,0,0),
  @(2.717167346176041,0.9067497345115783,0.2530916678862667,0,1.05543150756678,0.6405555441616428,0.0001344270136414139,0.0004405879691771197,0.4622344963235605,0.3179156297428563,0.4128635992005627,0,0,0),
  @(2.605915882671638,0.8714698395674816,0.2442637531186307,0,1.031141486411826,0.6255005086403429,0.0002738052848170192,0.0005747415003813927,0.4567224856003946,0.3175792512842044,0.3974470979293727,0,0,0),
  @(2.587410300903173,0.8671988759836324,0.2431604901138797,0,1.025357697903274,0.6213914440104844,0.0003025872047905942,0.000690337887167658,0.4549775764595836,0.3166709042556874,0.3948705170921585,0,0,0),
  @(2.715583298412696,0.9105952800541104,0.2539598757437886,0,1.050289926173463,0.6359243594328916,0.0001371955150284432,0.0006674638673382915,0.4598760355572011,0.3155671451993065,0.4126078882928397,0,0,0),
  @(3.281053413251925,1.099349943100435,0.3009635879206058,0,1.166524440641822,0.7058354673969376,0.0001573294586330043,0.00105600346457102,0.4846592228472986,0.3138011297402699,0.4910967277634768,0,0,0),
  @(4.391510873852951,1.468034860850537,0.3927468296867573,0,1.410511399277553,0.85789475611071,0.00327444172567759,0.004494834207204867,0.5415296330860997,0.3181310785706657,0.6463801873910029,0,0,0),
  @(5.209581941684007,1.739902649036821,0.4505898555455303,0,1.56047929074785,0.9528940717211754,0.007658998650318161,0.009271786826650086,0.5752837348679947,0.3152772092332512,0.7281919235798711,0,0,0),
  @(5.574915507912635,1.814655128295271,0.3848092003301531,0,1.346207466963293,0.8212084085000413,0.02595876379970008,0.01131445471829284,0.5030595373910671,0.2500099744106521,0.5033353314237132,0,0,0),
  @(5.710380277681679,1.811986786188811,0.3211623444235414,0,1.153852813542755,0.7042716506842197,0.06387214784956541,0.01156670974345531,0.4418435611642337,0.2047812531821549,0.3374157051991347,0,0,0),
  @(5.675201859578124,1.756535776282476,0.2574287918433953,0,0.960492820405392,0.5851708100276483,0.1184414230186093,0.01070523528458267,0.3820080749172092,0.1684783606773657,0.2079011622919396,0,0,0),
  @(5.573381726655157,1.695497586467809,0.2141844076723345,0,0.8266947849981818,0.5016603021039003,0.166858865779318,0.009725169245131227,0.3410679397424161,0.1475963609532087,0.1399196179894417,0,0,0),
  @(5.512175343742285,1.671382256161223,0.2032651434353738,0,0.7912694706717502,0.4788700082788466,0.1790075399511863,0.009365182933352756,0.3304034593838168,0.1432383157190209,0.1259694206843101,0,0,0),
  @(5.1658763999597,1.569286126986071,0.193809344595806,0,0.7680281498036408,0.4618560197332755,0.1647757305327531,0.007452901237233966,0.3274989131366368,0.1505316585761687,0.1221356014554829,0,0,0),
  @(4.955059478247904,1.523619909075933,0.2087365489762618,0,0.8191520131586287,0.491984243951137,0.1263380149986091,0.006567513224891464,0.3465266321035472,0.1666641996965144,0.1507928917996466,0,0,0),
  @(4.836077256499323,1.51701736364663,0.248740231843243,0,0.9499705868058754,0.572831057274982,0.07353491356904129,0.006092321877863682,0.3898487186640835,0.1962293063478882,0.2309928704082367,0,0,0),
  @(4.799141797284904,1.550574086639642,0.3123846854226002,0,1.140215457309452,0.6890884644148656,0.02907622708252688,0.006465108717167212,0.4494292494870393,0.2362402120355895,0.3780235155609688,0,0,0),
  @(4.994077813737533,1.681768864628907,0.4379940966766469,0,1.504299128693376,0.9129107280281232,0.006301994021111756,0.008500286326979634,0.5587253419593452,0.3084262609908244,0.7050615765017767,0,0,0),
  @(5.624909681759675,1.899710875596497,0.499693646957212,0,1.676035237888811,1.02228869653942,0.01073331568528424,0.01303123527801109,0.6029992529659154,0.3205703882858089,0.8206198861750948,0,0,0),
  @(6.038311166644178,2.034474411213125,0.5332721122954354,0,1.77951004149233,1.090846225474877,0.01409569977614467,0.01624200886630423,0.6305639031710086,0.3282985403205245,0.8795687504875787,0,0,0),
  @(5.817679291560864,1.957104655848184,0.5141103617270346,0,1.729833014986582,1.05934762797925,0.01226011426200119,0.01422913720236352,0.6184460957355071,0.3267983494599918,0.8481175364469351,0,0,0),
  @(4.984689579166343,1.675284041958662,0.4441039294689801,0,1.536290988626718,0.9351095708247072,0.006375273213042654,0.008101916767744299,0.5700846520420839,0.3183438434969368,0.7299172345772007,0,0,0),
  @(4.090529833339701,1.375459492816333,0.3695523648357835,0,1.334813500798674,0.8078327168787212,0.00205715680788332,0.003565944054884973,0.5215178567343628,0.312343339203764,0.604078579631576,0,0,0)
)

for ($r = 0; $r -lt 24; $r++) {
    $rowData = $data[$r]
    for ($c = 0; $c -lt 14; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowData[$c]
    }
    $ws.Cells.Item($r + 2, 16).Value = 0
    $ws.Cells.Item($r + 2, 17).Value = 0
}
